$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.01521177410757085
$ws.Range("C2").Value = 0.9708167685734876
$ws.Range("D2").Value = 4.187158799028484
$ws.Range("E2").Value = 2.046254822603598
$ws.Range("F2").Value = 2.066558963343684
$ws.Range("G2").Value = 51

# Row 3 (Q1)
$ws.Range("B3").Value = 0.06844627047333725
$ws.Range("C3").Value = 1.092801458624372
$ws.Range("D3").Value = 4.503710726994156
$ws.Range("E3").Value = 2.122194790068564
$ws.Range("F3").Value = 2.142625186267777
$ws.Range("G3").Value = 50

# Row 4 (Q2)
$ws.Range("B4").Value = 0.02116772330918306
$ws.Range("C4").Value = 0.973329744212709
$ws.Range("D4").Value = 3.945249452684095
$ws.Range("E4").Value = 1.986265202002012
$ws.Range("F4").Value = 2.006734846027487
$ws.Range("G4").Value = 49

# Row 5 (Q3)
$ws.Range("B5").Value = 0.10050846687517
$ws.Range("C5").Value = 1.118427419803017
$ws.Range("D5").Value = 4.727028183880057
$ws.Range("E5").Value = 2.17417298849012
$ws.Range("F5").Value = 2.194831739464666
$ws.Range("G5").Value = 48

# Row 6 (Q4)
$ws.Range("B6").Value = 0.0001922921193527008
$ws.Range("C6").Value = 0.9715821074844706
$ws.Range("D6").Value = 4.036252245157337
$ws.Range("E6").Value = 2.009042619049516
$ws.Range("F6").Value = 2.030762620643389
$ws.Range("G6").Value = 47

# Row 7 (Q5)
$ws.Range("B7").Value = 0.1002765460674019
$ws.Range("C7").Value = 1.08588965278254
$ws.Range("D7").Value = 4.674164069693668
$ws.Range("E7").Value = 2.161981514651239
$ws.Range("F7").Value = 2.183519073341691
$ws.Range("G7").Value = 46

# Row 8 (Q6)
$ws.Range("B8").Value = -0.02466186912181194
$ws.Range("C8").Value = 0.9115712867521671
$ws.Range("D8").Value = 4.026575163211482
$ws.Range("E8").Value = 2.006632792319382
$ws.Range("F8").Value = 2.029154061280164
$ws.Range("G8").Value = 45
